$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 66.85111099999999
$ws.Range("H2").Value = 200.553333
$ws.Range("I2").Value = 0.1215550702639512
$ws.Range("J2").Value = 0.1215550702639512
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 73.19179533333333
$ws.Range("N2").Value = 219.575386
$ws.Range("O2").Value = 0.4454729128883617
$ws.Range("P2").Value = 0.4454729128883617
$ws.Range("Q2").Value = 4892.952834117948
$ws.Range("R2").Value = 44036.57550706153
$ws.Range("S2").Value = 0.05414949122683183
$ws.Range("T2").Value = 0.05414949122683183
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 66.85111099999999
$ws.Range("H3").Value = 200.553333
$ws.Range("I3").Value = 0.1215550702639512
$ws.Range("J3").Value = 0.1215550702639512
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.224257666666666
$ws.Range("N3").Value = 21.672773
$ws.Range("O3").Value = 0.04396956095378667
$ws.Range("P3").Value = 0.04396956095378668
$ws.Range("Q3").Value = 482.9496511669342
$ws.Range("R3").Value = 4346.546860502409
$ws.Range("S3").Value = 0.005344723071212625
$ws.Range("T3").Value = 0.005344723071212626
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 66.85111099999999
$ws.Range("H4").Value = 200.553333
$ws.Range("I4").Value = 0.1215550702639512
$ws.Range("J4").Value = 0.1215550702639512
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 81.09049733333333
$ws.Range("N4").Value = 243.271492
$ws.Range("O4").Value = 0.4935473968103956
$ws.Range("P4").Value = 0.4935473968103956
$ws.Range("Q4").Value = 5420.989838275869
$ws.Range("R4").Value = 48788.90854448283
$ws.Range("S4").Value = 0.05999318849787785
$ws.Range("T4").Value = 0.05999318849787785
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 66.85111099999999
$ws.Range("H5").Value = 200.553333
$ws.Range("I5").Value = 0.1215550702639512
$ws.Range("J5").Value = 0.1215550702639512
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.794787
$ws.Range("N5").Value = 8.384361
$ws.Range("O5").Value = 0.01701012934745599
$ws.Range("P5").Value = 0.01701012934745599
$ws.Range("Q5").Value = 186.8346159583569
$ws.Range("R5").Value = 1681.511543625213
$ws.Range("S5").Value = 0.002067667468028912
$ws.Range("T5").Value = 0.002067667468028912
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 99.86393233333332
$ws.Range("H6").Value = 299.591797
$ws.Range("I6").Value = 0.1815821327429069
$ws.Range("J6").Value = 0.1815821327429069
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 73.19179533333333
$ws.Range("N6").Value = 219.575386
$ws.Range("O6").Value = 0.4454729128883617
$ws.Range("P6").Value = 0.4454729128883617
$ws.Range("Q6").Value = 7309.220496523181
$ws.Range("R6").Value = 65782.98446870863
$ws.Range("S6").Value = 0.0808899216014639
$ws.Range("T6").Value = 0.08088992160146391
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 99.86393233333332
$ws.Range("H7").Value = 299.591797
$ws.Range("I7").Value = 0.1815821327429069
$ws.Range("J7").Value = 0.1815821327429069
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.224257666666666
$ws.Range("N7").Value = 21.672773
$ws.Range("O7").Value = 0.04396956095378667
$ws.Range("P7").Value = 0.04396956095378668
$ws.Range("Q7").Value = 721.4427787825645
$ws.Range("R7").Value = 6492.985009043081
$ws.Range("S7").Value = 0.007984086653757829
$ws.Range("T7").Value = 0.007984086653757829
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 99.86393233333332
$ws.Range("H8").Value = 299.591797
$ws.Range("I8").Value = 0.1815821327429069
$ws.Range("J8").Value = 0.1815821327429069
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 81.09049733333333
$ws.Range("N8").Value = 243.271492
$ws.Range("O8").Value = 0.4935473968103956
$ws.Range("P8").Value = 0.4935473968103956
$ws.Range("Q8").Value = 8098.015938572346
$ws.Range("R8").Value = 72882.14344715112
$ws.Range("S8").Value = 0.0896193889225414
$ws.Range("T8").Value = 0.08961938892254141
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 99.86393233333332
$ws.Range("H9").Value = 299.591797
$ws.Range("I9").Value = 0.1815821327429069
$ws.Range("J9").Value = 0.1815821327429069
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.794787
$ws.Range("N9").Value = 8.384361
$ws.Range("O9").Value = 0.01701012934745599
$ws.Range("P9").Value = 0.01701012934745599
$ws.Range("Q9").Value = 279.0984198540796
$ws.Range("R9").Value = 2511.885778686717
$ws.Range("S9").Value = 0.00308873556514377
$ws.Range("T9").Value = 0.003088735565143771
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 368.4456226666667
$ws.Range("H10").Value = 1105.336868
$ws.Range("I10").Value = 0.6699429954379058
$ws.Range("J10").Value = 0.6699429954379058
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 73.19179533333333
$ws.Range("N10").Value = 219.575386
$ws.Range("O10").Value = 0.4454729128883617
$ws.Range("P10").Value = 0.4454729128883617
$ws.Range("Q10").Value = 26967.19660568123
$ws.Range("R10").Value = 242704.7694511311
$ws.Range("S10").Value = 0.2984414576468783
$ws.Range("T10").Value = 0.2984414576468783
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 368.4456226666667
$ws.Range("H11").Value = 1105.336868
$ws.Range("I11").Value = 0.6699429954379058
$ws.Range("J11").Value = 0.6699429954379058
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 7.224257666666666
$ws.Range("N11").Value = 21.672773
$ws.Range("O11").Value = 0.04396956095378667
$ws.Range("P11").Value = 0.04396956095378668
$ws.Range("Q11").Value = 2661.74611429944
$ws.Range("R11").Value = 23955.71502869497
$ws.Range("S11").Value = 0.02945709937346943
$ws.Range("T11").Value = 0.02945709937346943
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 368.4456226666667
$ws.Range("H12").Value = 1105.336868
$ws.Range("I12").Value = 0.6699429954379058
$ws.Range("J12").Value = 0.6699429954379058
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 81.09049733333333
$ws.Range("N12").Value = 243.271492
$ws.Range("O12").Value = 0.4935473968103956
$ws.Range("P12").Value = 0.4935473968103956
$ws.Range("Q12").Value = 29877.43878232967
$ws.Range("R12").Value = 268896.9490409671
$ws.Range("S12").Value = 0.3306486214097371
$ws.Range("T12").Value = 0.3306486214097372
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 368.4456226666667
$ws.Range("H13").Value = 1105.336868
$ws.Range("I13").Value = 0.6699429954379058
$ws.Range("J13").Value = 0.6699429954379058
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.794787
$ws.Range("N13").Value = 8.384361
$ws.Range("O13").Value = 0.01701012934745599
$ws.Range("P13").Value = 0.01701012934745599
$ws.Range("Q13").Value = 1029.727036435705
$ws.Range("R13").Value = 9267.543327921348
$ws.Range("S13").Value = 0.0113958170078209
$ws.Range("T13").Value = 0.0113958170078209
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 14.80496566666667
$ws.Range("H14").Value = 44.414897
$ws.Range("I14").Value = 0.02691980155523597
$ws.Range("J14").Value = 0.02691980155523597
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 73.19179533333333
$ws.Range("N14").Value = 219.575386
$ws.Range("O14").Value = 0.4454729128883617
$ws.Range("P14").Value = 0.4454729128883617
$ws.Range("Q14").Value = 1083.602016991694
$ws.Range("R14").Value = 9752.418152925242
$ws.Range("S14").Value = 0.01199204241318762
$ws.Range("T14").Value = 0.01199204241318762
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 14.80496566666667
$ws.Range("H15").Value = 44.414897
$ws.Range("I15").Value = 0.02691980155523597
$ws.Range("J15").Value = 0.02691980155523597
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 7.224257666666666
$ws.Range("N15").Value = 21.672773
$ws.Range("O15").Value = 0.04396956095378667
$ws.Range("P15").Value = 0.04396956095378668
$ws.Range("Q15").Value = 106.9548867221534
$ws.Range("R15").Value = 962.593980499381
$ws.Range("S15").Value = 0.001183651855346789
$ws.Range("T15").Value = 0.001183651855346789
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 14.80496566666667
$ws.Range("H16").Value = 44.414897
$ws.Range("I16").Value = 0.02691980155523597
$ws.Range("J16").Value = 0.02691980155523597
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 81.09049733333333
$ws.Range("N16").Value = 243.271492
$ws.Range("O16").Value = 0.4935473968103956
$ws.Range("P16").Value = 0.4935473968103956
$ws.Range("Q16").Value = 1200.542028912925
$ws.Range("R16").Value = 10804.87826021632
$ws.Range("S16").Value = 0.01328619798023915
$ws.Range("T16").Value = 0.01328619798023915
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 14.80496566666667
$ws.Range("H17").Value = 44.414897
$ws.Range("I17").Value = 0.02691980155523597
$ws.Range("J17").Value = 0.02691980155523597
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.794787
$ws.Range("N17").Value = 8.384361
$ws.Range("O17").Value = 0.01701012934745599
$ws.Range("P17").Value = 0.01701012934745599
$ws.Range("Q17").Value = 41.37672558064634
$ws.Range("R17").Value = 372.390530225817
$ws.Range("S17").Value = 0.0004579093064624108
$ws.Range("T17").Value = 0.0004579093064624108
